# Apply the "added ability to leave sites out of backend upgrades, added
# cost example for srr" edit to the SiteDevelopmentValues sheet.
#
# Net visible changes (per the OOXML diff):
#   - receiver_cost_factor row (13): cost guess -> value sourced from the
#     "station cost estimate" sheet (Kari); value 1,500,000 -> 1,300,000
#   - maser_cost row (16) renamed to timing_dbe_cost, description stays
#     "Maser"; value 300,000 -> 850,000; note updated to the same
#     "station cost estimate" source
#   - selection cursor left on A16

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SiteDevelopmentValues")

# Row 13: receiver_cost_factor
$ws.Range("C13").Value = 1300000
$ws.Range("D13").Value = 'from "station cost estimate" sheet (Kari)'

# Row 16: maser_cost -> timing_dbe_cost
$ws.Range("A16").Value = "timing_dbe_cost"
$ws.Range("C16").Value = 850000
$ws.Range("D16").Value = 'from "station cost estimate" sheet (Kari)'

# Move / leave the active selection on A16, matching the saved view state.
$ws.Range("A16").Select() | Out-Null
